# Applies the "display player's suns using SunImageView" edit to the
# "Created custom ImageView class to drawText a number on top of a Sun
# Tile" bullet item.
#
# Strategy: locate the target paragraph by its (unique) text, replace its
# inner run content with the expanded wording ("Sun" + "ImageView" new
# wording, "(extends ImageView)" aside), then insert a brand-new list
# paragraph right after it describing SunImageView's use for displaying
# player's suns. Finally relocate the (hidden) _GoBack bookmark to sit at
# the end of the newly inserted sentence, mirroring where Word leaves it
# after the last edit.

$d = $word.ActiveDocument

# --- locate the paragraph that needs editing -------------------------------
# NB: Paragraph.Range.Text includes the trailing paragraph mark (CR, \r) --
# always trim it before doing exact text comparisons.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Created custom ImageView class to drawText a number on top of a Sun Tile") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Created custom ImageView...' paragraph"
}

# --- detach the _GoBack bookmark so it doesn't get dragged around ----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$flatOpcHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$flatOpcFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- rewrite the body of the first paragraph --------------------------------
$newFirstParaRuns =
    '<w:r><w:t xml:space="preserve">Created custom </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Sun</w:t></w:r>' +
    '<w:r><w:t>ImageView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">(extends </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ImageView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">) </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">class to </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>drawText</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> a number on top of a Sun Tile</w:t></w:r>'

$firstParaXml = $flatOpcHeader + '<w:body><w:p>' + $newFirstParaRuns + '</w:p></w:body>' + $flatOpcFooter

$innerRange = $d.Range($target.Range.Start, $target.Range.End)
$innerRange.InsertXML($firstParaXml)

# --- insert the brand-new list paragraph right after it ---------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Created custom SunImageView (extends ImageView) class to drawText a number on top of a Sun Tile") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not re-locate the rewritten paragraph"
}

$newSecondParaXml = $flatOpcHeader +
    '<w:body><w:p>' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Use </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>SunImageView</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to display player' + [char]8217 + 's suns and sun in auction.</w:t></w:r>' +
    '</w:p></w:body>' + $flatOpcFooter

$afterRange = $d.Range($target.Range.End, $target.Range.End)
$afterRange.InsertXML($newSecondParaXml)

# --- find the freshly inserted paragraph and re-anchor _GoBack at its end --
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq ("Use SunImageView to display player" + [char]8217 + "s suns and sun in auction.")) {
        $newPara = $p
        break
    }
}
if ($newPara -eq $null) {
    throw "Could not locate the newly inserted SunImageView paragraph"
}

$goBackPos = $newPara.Range.End
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
